$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Recommandations sheet: updated metrics / reshuffled rows (rows 2-45, 47) ---
# Row 2
$ws1.Range("D2").Value = 3329.06
$ws1.Range("E2").Value = 104.39
# Row 3
$ws1.Range("D3").Value = 2755.0
$ws1.Range("E3").Value = 680.0
# Row 4
$ws1.Range("D4").Value = 2660.0
$ws1.Range("E4").Value = 660.0
# Row 5
$ws1.Range("D5").Value = 2479.45
$ws1.Range("E5").Value = 571.79
# Row 6
$ws1.Range("D6").Value = 1710.79
$ws1.Range("E6").Value = 422.5
# Row 7
$ws1.Range("D7").Value = 1458.84
$ws1.Range("E7").Value = 364.71
# Row 8
$ws1.Range("D8").Value = 1444.51
$ws1.Range("E8").Value = 354.06
# Row 9
$ws1.Range("D9").Value = 557.64
$ws1.Range("E9").Value = 134.23
# Row 10
$ws1.Range("D10").Value = 546.93
$ws1.Range("E10").Value = 136.55
# Row 11
$ws1.Range("D11").Value = 541.83
$ws1.Range("E11").Value = 135.92
# Row 12
$ws1.Range("D12").Value = 532.5
$ws1.Range("E12").Value = 133.58
# Row 13
$ws1.Range("D13").Value = 493.99
$ws1.Range("E13").Value = 121.66
# Row 14
$ws1.Range("D14").Value = 429.94
$ws1.Range("E14").Value = 107.29
# Row 15
$ws1.Range("D15").Value = 379.62
$ws1.Range("E15").Value = 94.91
# Row 16
$ws1.Range("A16").Value = "BRVM - INDUSTRIE                  (**)"
$ws1.Range("D16").Value = 207.08
$ws1.Range("E16").Value = 207.08
# Row 19
$ws1.Range("A19").Value = "BRVM-PRINCIPAL                    (**)"
$ws1.Range("D19").Value = 193.4
$ws1.Range("E19").Value = 193.4
# Row 21
$ws1.Range("A21").Value = "BRVM - CONSOMMATION DE BASE         (**)"
$ws1.Range("D21").Value = 185.56
$ws1.Range("E21").Value = 185.56
# Row 22
$ws1.Range("D22").Value = 18.76
$ws1.Range("E22").Value = 7.41
# Row 23
$ws1.Range("A23").Value = "FILTISAC CI (FTSC)"
$ws1.Range("B23").Value = 3.0
$ws1.Range("C23").Value = 1.0
$ws1.Range("D23").Value = 14.76
$ws1.Range("E23").Value = 7.46
$ws1.Range("F23").Value = "🟢 Achat"
$ws1.Range("G23").Value = "✅ Renforcer"
# Row 25
$ws1.Range("A25").Value = "SETAO CI (STAC)"
$ws1.Range("D25").Value = 5.42
$ws1.Range("E25").Value = 5.42
# Row 26
$ws1.Range("A26").Value = "BANK OF AFRICA BF (BOABF)"
$ws1.Range("B26").Value = 1.0
$ws1.Range("D26").Value = 4.85
$ws1.Range("E26").Value = 4.85
# Row 28
$ws1.Range("A28").Value = "SODE CI (SDCC)"
$ws1.Range("C28").Value = 1.0
$ws1.Range("D28").Value = 2.9
$ws1.Range("E28").Value = 4.35
$ws1.Range("G28").Value = "👀 À surveiller"
# Row 29
$ws1.Range("A29").Value = "ONATEL BF (ONTBF)"
$ws1.Range("C29").Value = 0.0
$ws1.Range("D29").Value = 2.61
$ws1.Range("E29").Value = 2.61
$ws1.Range("G29").Value = "➖ Neutre"
# Row 30
$ws1.Range("A30").Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws1.Range("C30").Value = 0.0
$ws1.Range("D30").Value = 2.5
$ws1.Range("E30").Value = 2.5
$ws1.Range("G30").Value = "➖ Neutre"
# Row 31
$ws1.Range("A31").Value = "BICI CI (BICC)"
$ws1.Range("D31").Value = 2.2
$ws1.Range("E31").Value = 2.2
# Row 32
$ws1.Range("A32").Value = "VIVO ENERGY CI (SHEC)"
$ws1.Range("B32").Value = 1.0
$ws1.Range("C32").Value = 0.0
$ws1.Range("D32").Value = 1.27
$ws1.Range("E32").Value = 1.27
# Row 33
$ws1.Range("A33").Value = "SAPH CI (SPHC)"
$ws1.Range("D33").Value = 1.22
$ws1.Range("E33").Value = -6.26
# Row 34
$ws1.Range("A34").Value = "SICABLE CI (CABC)"
$ws1.Range("B34").Value = 1.0
$ws1.Range("C34").Value = 1.0
$ws1.Range("D34").Value = 0.25
$ws1.Range("E34").Value = 5.2
# Row 35
$ws1.Range("A35").Value = "TOTAL"
$ws1.Range("B35").Value = 0.0
$ws1.Range("C35").Value = 3.0
$ws1.Range("D35").Value = 0.0
$ws1.Range("E35").Value = 0.0
$ws1.Range("G35").Value = "➖ Neutre"
# Row 36
$ws1.Range("A36").Value = "BERNABE CI (BNBC)"
$ws1.Range("B36").Value = 1.0
$ws1.Range("D36").Value = -1.01
$ws1.Range("E36").Value = 6.45
$ws1.Range("G36").Value = "👀 À surveiller"
# Row 37
$ws1.Range("A37").Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$ws1.Range("D37").Value = -1.5
$ws1.Range("E37").Value = -1.5
# Row 38
$ws1.Range("A38").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws1.Range("B38").Value = 0.0
$ws1.Range("D38").Value = -1.73
$ws1.Range("E38").Value = -1.73
$ws1.Range("G38").Value = "➖ Neutre"
# Row 39
$ws1.Range("A39").Value = "CFAO MOTORS CI (CFAC)"
$ws1.Range("B39").Value = 0.0
$ws1.Range("D39").Value = -1.78
$ws1.Range("E39").Value = -1.78
$ws1.Range("G39").Value = "➖ Neutre"
# Row 40
$ws1.Range("A40").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Range("D40").Value = -2.04
$ws1.Range("E40").Value = -2.04
# Row 41
$ws1.Range("A41").Value = "UNILEVER CI (UNLC)"
$ws1.Range("B41").Value = 1.0
$ws1.Range("D41").Value = -2.24
$ws1.Range("E41").Value = -7.5
$ws1.Range("G41").Value = "👀 À surveiller"
# Row 42
$ws1.Range("A42").Value = "SICOR CI (SICC)"
$ws1.Range("B42").Value = 1.0
$ws1.Range("D42").Value = -2.47
$ws1.Range("E42").Value = -6.91
$ws1.Range("G42").Value = "👀 À surveiller"
# Row 43
$ws1.Range("A43").Value = "CIE CI (CIEC)"
$ws1.Range("B43").Value = 0.0
$ws1.Range("C43").Value = 1.0
$ws1.Range("D43").Value = -2.54
$ws1.Range("E43").Value = -2.54
$ws1.Range("G43").Value = "➖ Neutre"
# Row 44
$ws1.Range("A44").Value = "UNIWAX CI (UNXC)"
$ws1.Range("D44").Value = -3.31
$ws1.Range("E44").Value = -3.31
# Row 45
$ws1.Range("A45").Value = "SUCRIVOIRE (SCRC)"
$ws1.Range("C45").Value = 2.0
$ws1.Range("D45").Value = -10.26
$ws1.Range("E45").Value = -2.92
# Row 47
$ws1.Range("A47").Value = "LOTERIE NATIONALE DU BENIN (LNBB)"
$ws1.Range("B47").Value = 0.0
$ws1.Range("C47").Value = 2.0
$ws1.Range("D47").Value = -14.92
$ws1.Range("E47").Value = -7.45
$ws1.Range("F47").Value = "🟡 Observer"
$ws1.Range("G47").Value = "➖ Neutre"

# --- Top_YTD sheet: updated YTD progression values (column B, rows 2-11) ---
$ws2.Range("B2").Value = 8243775.08
$ws2.Range("B3").Value = 386873.6
$ws2.Range("B4").Value = 342286.46
$ws2.Range("B5").Value = 267354.51
$ws2.Range("B6").Value = 77433.08
$ws2.Range("B7").Value = 46535.43
$ws2.Range("B8").Value = 45105.36
$ws2.Range("B9").Value = 3183.6
$ws2.Range("B10").Value = 3040.72
$ws2.Range("B11").Value = 2973.61
